$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 22:35"

# Update numeric data for countries whose totals changed
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1662131
$ws.Range("C4").Value = 17037
$ws.Range("D4").Value = 444381
$ws.Range("E4").Value = 1119190
$ws.Range("G4").Value = 913
$ws.Range("H4").Value = 98560

# Peru (row 15)
$ws.Range("B15").Value = 115754
$ws.Range("C15").Value = 4056
$ws.Range("D15").Value = 47915
$ws.Range("E15").Value = 64466
$ws.Range("G15").Value = 129
$ws.Range("H15").Value = 3373

# Canada (row 16)
$ws.Range("B16").Value = 83593
$ws.Range("C16").Value = 1113
$ws.Range("D16").Value = 43222
$ws.Range("E16").Value = 34019

# Niger (row 113)
$ws.Range("B113").Value = 943
$ws.Range("C113").Value = 6
$ws.Range("D113").Value = 775
$ws.Range("E113").Value = 107
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 61

# Estado de Palestina (row 139)
$ws.Range("D139").Value = 348
$ws.Range("E139").Value = 72

# Monaco overtakes Bahamas in total cases, so they swap rows (170/171)
# Row 170 becomes Monaco, row 171 becomes Bahamas
$ws.Range("A170").Value = "Monaco"
$ws.Range("B170").Value = 98
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 90
$ws.Range("E170").Value = 4
$ws.Range("H170").Value = 4

$ws.Range("A171").Value = "Bahamas"
$ws.Range("D171").Value = 45
$ws.Range("E171").Value = 41
$ws.Range("H171").Value = 11
